# Apply "Add files via upload" changes:
# - Update Year (column A) from 2021 to 2022 for rows 126-137
# - Add DraftPosition (column H) and Bye (column I) values for rows 126-137
# - Update the active sheet view (topLeftCell / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Draft position (H) and Bye (I) values per row, keyed by row number
$rowData = @{
    126 = @{ H = 12; I = 0 }
    127 = @{ H = 7;  I = 0 }
    128 = @{ H = 2;  I = 0 }
    129 = @{ H = 5;  I = 0 }
    130 = @{ H = 9;  I = 0 }
    131 = @{ H = 11; I = 0 }
    132 = @{ H = 3;  I = 0 }
    133 = @{ H = 10; I = 0 }
    134 = @{ H = 6;  I = 0 }
    135 = @{ H = 4;  I = 1 }
    136 = @{ H = 8;  I = 1 }
    137 = @{ H = 1;  I = 0 }
}

foreach ($r in 126..137) {
    # Year column A: 2021 -> 2022
    $ws.Cells.Item($r, 1).Value = 2022

    # New columns H (DraftPosition) and I (Bye) - match the centered style
    # already used throughout the rest of the sheet (e.g. column G).
    $data = $rowData[$r]

    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $data.H
    $hCell.HorizontalAlignment = $ws.Cells.Item($r, 7).HorizontalAlignment

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.Value = $data.I
    $iCell.HorizontalAlignment = $ws.Cells.Item($r, 7).HorizontalAlignment
}

# Update sheet view: topLeftCell moves from A115 to B115, selection moves to I135
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 115
$win.ScrollColumn = 2
$ws.Range("I135").Select()

# Best-effort: workbook window geometry also shifted in the saved file
# (xWindow 11520 -> 0, windowWidth 11520 -> 13560).
$win.Left = 0
$win.Width = 13560
